# BIS-769: Fixed xls test files
# Add "Pattern" / "Pattern Type" header columns (M4, N4) to the
# SAMPLE_TYPE export sheet, matching the style already used by the
# other header cells (K4 / L4), and rewrite the "Unique" boolean
# column (L5:L7) as a FALSE() formula instead of a literal boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Match the formatting of the neighbouring header cell (K4/L4)
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)

# The "Unique" column values become FALSE() formulas
$ws.Range("L5").Formula = "=FALSE()"
$ws.Range("L6").Formula = "=FALSE()"
$ws.Range("L7").Formula = "=FALSE()"

# Update the active selection to the new header cells
$ws.Range("M4:N4").Select()
